$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# New data rows 2-7 for the "dang ky chuc nang" (function registration)
# table. Column order: A=Ten nhom, B=Thoi gian hoan thanh, C=Chuc nang,
# D=Thanh vien hoac nhom phat trien.
# ------------------------------------------------------------------

# Row 6/7's "Thoi gian hoan thanh" column holds a real date value
# (centered, mm-dd-yy number format -> numFmtId 14) rather than literal
# text like rows 2-5. Build that style once on B6 and clone it onto B7
# via Copy/PasteSpecial(xlPasteFormats) so both cells share a single
# cellXfs entry (matches how Excel itself interns identical styles).
$ws.Range("B6").HorizontalAlignment = -4108
$ws.Range("B6").NumberFormat = "mm-dd-yy"
$ws.Range("B6").Value = 43090

$ws.Range("B6").Copy()
$ws.Range("B7").PasteSpecial(-4122)
$ws.Range("B7").Value = 43090

# Rows 2-5 store the completion date as plain centered text.
$ws.Range("B2:B5").HorizontalAlignment = -4108

$ws.Range("A2").Value = "The Owls"
$ws.Range("B2").Value = "20/12/2017"
$ws.Range("C2").Value = "Quản lí nhân viên(xem, xóa, sửa, thêm)"
$ws.Range("D2").Value = "Nguyễn Hoàng Quân(1412439)"

$ws.Range("A3").Value = "The Owls"
$ws.Range("B3").Value = "21/12/2017"
$ws.Range("C3").Value = "Quản lý khách hàng(xem, sửa, thêm)"
$ws.Range("D3").Value = "Ngô Thị Mai Lý (1412310)"

$ws.Range("A4").Value = "The Owls"
$ws.Range("B4").Value = "21/12/2017"
$ws.Range("C4").Value = "Quản lý thống kê"
$ws.Range("D4").Value = "Võ Ngọc Công Minh (1412328)"

$ws.Range("A5").Value = "The Owls"
$ws.Range("B5").Value = "21/12/2017"
$ws.Range("C5").Value = "Quán lý thuê phòng (xem, sửa, thêm)"
$ws.Range("D5").Value = "Trần Nguyên (1412360)"

$ws.Range("A6").Value = "The Owls"
$ws.Range("C6").Value = "Quản lý dịch vụ (xem, sửa, thêm)"
$ws.Range("D6").Value = "Trần Trọng Cao Nguyên (1412359)"

$ws.Range("A7").Value = "The Owls"
$ws.Range("C7").Value = "Quản lý thiết bị(xem, thêm, xóa, sửa)"
$ws.Range("D7").Value = "Võ Đặng Nguyễn(1412362)"

# Widen the "Chuc nang" column (C) from 17.71 to 40 characters.
$ws.Range("C1").ColumnWidth = 39.1666666666667

# Selection ends up parked one row below the new data, like the source file.
$ws.Range("D8").Select() | Out-Null
